$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B from 10 to 11 (COM ColumnWidth is offset by -0.83 from the
# OOXML <col width> value for this workbook's default font).
$ws.Columns.Item(2).ColumnWidth = 10.17

# Append the new ranking row at the bottom of the table (row 359).
$ws.Cells.Item(359, 1).Value = "2025/12/18 05:00"
$ws.Cells.Item(359, 2).Value = "109,230位本"
$ws.Cells.Item(359, 3).Value = "192位 広告・宣伝 (本)"
$ws.Cells.Item(359, 4).Value = "327位商業デザイン"
$ws.Cells.Item(359, 5).Value = "4,159位ビジネス実用本"
$ws.Cells.Item(359, 6).Value = "-"
$ws.Cells.Item(359, 7).Value = "-"
